$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("BU") values were replaced as part of switching the data storage
# to RDS format. Values are stored as text (not numbers), so force a text
# number format before assigning the new values to preserve that type.
$ws.Range("D2:D9").NumberFormat = "@"

$ws.Range("D2").Value = "220"
$ws.Range("D3").Value = "53"
$ws.Range("D4").Value = "99"
$ws.Range("D5").Value = "131"
$ws.Range("D6").Value = "51"
$ws.Range("D7").Value = "35"
$ws.Range("D8").Value = "68"
$ws.Range("D9").Value = "657"
